$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "256.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.39%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.84"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-0.78%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.530"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.91%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.42%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.610"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-1.01%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8508"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.24%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9289"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.72%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.00%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.04518"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "25.07%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07088"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.02%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03067"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.17%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09086"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.66%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001537"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.75%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006029"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.33%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006103"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.81%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.481"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.20%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.170"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.24%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3049"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.64%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.58%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.900"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.24%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04254"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.75%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.14%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004303"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.36%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.05%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "2.01%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03801"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.46%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006228"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.92%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1099"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.18%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002429"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.85%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "31.03%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005350"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.74%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.03%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-50.62%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2525"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "10,992.56%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.03%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.03%"
